# Weekly update: insert two new rows of price data at the top of the
# "Vega Modelo de Temuco - Cebollín" data block (rows 583-584), pushing
# the existing data down by two rows. This mirrors a new week's data
# being prepended to the consolidated sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 583 (shifts old 583..666
# down to 585..668, growing the used range to A1:R668).
$ws.Rows("583:584").Insert()

# --- New row 583 ---------------------------------------------------
# Same market/category metadata as the (now shifted) row below it, but
# with the new week's date and a new volume figure.
$ws.Range("A583").Value = 10
$ws.Range("B583").Value = "Vega Modelo de Temuco"
$ws.Range("C583").Value = "La Araucanía"
$ws.Range("D583").Value = 45131
$ws.Range("E583").Value = 9
$ws.Range("F583").Value = 100112037
$ws.Range("G583").Value = "Cebollín"
$ws.Range("H583").Value = "Sin especificar"
$ws.Range("I583").Value = "Primera"
$ws.Range("J583").Value = 125
$ws.Range("K583").Value = 8000
$ws.Range("L583").Value = 8000
$ws.Range("M583").Value = 8000
$ws.Range("N583").Value = "$/docena de paquetes"
$ws.Range("O583").Value = "Provincia de Cautín"
$ws.Range("P583").Value = 667
$ws.Range("Q583").Value = 12
$ws.Range("R583").Value = "Hortaliza"

# --- New row 584 ---------------------------------------------------
$ws.Range("A584").Value = 10
$ws.Range("B584").Value = "Vega Modelo de Temuco"
$ws.Range("C584").Value = "La Araucanía"
$ws.Range("D584").Value = 45131
$ws.Range("E584").Value = 9
$ws.Range("F584").Value = 100112037
$ws.Range("G584").Value = "Cebollín"
$ws.Range("H584").Value = "Sin especificar"
$ws.Range("I584").Value = "Primera"
$ws.Range("J584").Value = 150
$ws.Range("K584").Value = 7000
$ws.Range("L584").Value = 7000
$ws.Range("M584").Value = 7000
$ws.Range("N584").Value = "$/docena de paquetes"
$ws.Range("O584").Value = "Región de O'Higgins"
$ws.Range("P584").Value = 583
$ws.Range("Q584").Value = 12
$ws.Range("R584").Value = "Hortaliza"
